$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add "I0" (column I) and "IF" (column J), matching the
#     existing header style (bold, centered, bordered -> same as H1). ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-35 ---
# Row 2 is a special case with hard-coded values.
$ws.Range("I2").Value = 13
$ws.Range("J2").Value = 15

# Remaining rows: I = 1 (constant) and J = copy of column H's value for
# that row.
for ($r = 3; $r -le 35; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}

Write-Output "Added I0/IF columns"
